# Apply the "Import sitepreview\ips-pilgrimage" metadata/content refresh to
# the ConsentVerifiedBy StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata": update a handful of property values and insert a new
# "Jurisdiction" row right after "Contact".
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2.0.2"            # Version: 1.0.0 -> 2.0.2
$meta.Range("B6").Value = "active"           # Status: draft -> active
$meta.Range("B8").Value = "2025-02-04T19:22:12+00:00"  # Date

# Insert the new "Jurisdiction" row after row 10 ("Contact"), shifting the
# remaining rows (Description, Purpose, Copyright, FHIR Version, Kind, Type,
# Base Definition, Abstract, Derivation, Context) down by one.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# FHIR Version moved from row 14 to row 15 because of the insert above.
$meta.Range("B15").Value = "4.0.1"           # FHIR Version: 4.3.0 -> 4.0.1

# ---------------------------------------------------------------------------
# Sheet "Elements": a few corrections to the Extension definition table.
# ---------------------------------------------------------------------------
$elem = $wb.Worksheets.Item("Elements")

# Extension.id (row 3) Type(s) column: id -> string
$elem.Range("K3").Value = "string
"

# Extension.value[x] (row 6) Definition column: R4B -> R4 doc link
$elem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."

# Extension (row 2) Constraint(s) column: drop the "unless an empty
# Parameters resource ... or `$this is Parameters" clause from ele-1.
$elem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
